$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Sheet1" to "Report"
$ws.Name = "Report"

# Update membership counts (column B) for GASCO, IOS, IOWA, MOASC
$ws.Range("B3").Value = 3200
$ws.Range("B4").Value = 1200
$ws.Range("B5").Value = 1760
$ws.Range("B6").Value = 4000

# Update the active selection to match the saved view state
$ws.Range("E10").Select()
